$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 297.85715
$ws.Range("I6").Value = 296.33334
$ws.Range("K6").Value = 889.0000200000001
$ws.Range("M6").Value = -777.0000200000001
$ws.Range("H12").Value = 1088.4
$ws.Range("I12").Value = 999
$ws.Range("J12").Value = 1110.75
$ws.Range("K12").Value = 999
$ws.Range("L12").Value = 1110.75
$ws.Range("M12").Value = -829
$ws.Range("N12").Value = -1450.75
$ws.Range("H17").Value = 1284.6842
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 1300.5
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3901.5
$ws.Range("N17").Value = -4237.5
$ws.Range("H21").Value = 6166.6665
$ws.Range("J21").Value = 6166.6665
$ws.Range("L21").Value = 6166.6665
$ws.Range("N21").Value = -7102.6665
$ws.Range("H23").Value = 6166.6665
$ws.Range("J23").Value = 6166.6665
$ws.Range("L23").Value = 6166.6665
$ws.Range("N23").Value = -6634.6665
$ws.Range("H62").Value = 1100
$ws.Range("I62").Value = 1100
$ws.Range("K62").Value = 1100
$ws.Range("M62").Value = -476
$ws.Range("H65").Value = 1100
$ws.Range("I65").Value = 1100
$ws.Range("K65").Value = 5500
$ws.Range("M65").Value = -2380
$ws.Range("H86").Value = 73786380
$ws.Range("I86").Value = 160714990
$ws.Range("J86").Value = 6175254
$ws.Range("K86").Value = 160714990
$ws.Range("L86").Value = 6175254
$ws.Range("M86").Value = -160713867
$ws.Range("N86").Value = -6177500
$ws.Range("H89").Value = 73786380
$ws.Range("I89").Value = 160714990
$ws.Range("J89").Value = 6175254
$ws.Range("K89").Value = 803574950
$ws.Range("L89").Value = 30876270
$ws.Range("M89").Value = -803569334
$ws.Range("N89").Value = -30887502
$ws.Range("H107").Value = 41668744
$ws.Range("I107").Value = 15627337
$ws.Range("K107").Value = 15627337
$ws.Range("M107").Value = -15625417
$ws.Range("H116").Value = 31260300
$ws.Range("J116").Value = 12582
$ws.Range("L116").Value = 12582
$ws.Range("N116").Value = -19466
$ws.Range("H131").Value = 1511.409
$ws.Range("I131").Value = 1162.55
$ws.Range("K131").Value = 3487.65
$ws.Range("M131").Value = 1552.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2904903
$ws.Range("I32").Value = 3129584.5
$ws.Range("K32").Value = 3129584.5
$ws.Range("M32").Value = -3129297.5
$ws.Range("H56").Value = 99999.664
$ws.Range("J56").Value = 99999.664
$ws.Range("L56").Value = 99999.664
$ws.Range("N56").Value = -101483.664
$ws.Range("H110").Value = 18519630
$ws.Range("I110").Value = 1149.4286
$ws.Range("K110").Value = 1149.4286
$ws.Range("M110").Value = 895.5714
$ws.Range("H114").Value = 56340
$ws.Range("J114").Value = 56340
$ws.Range("L114").Value = 56340
$ws.Range("N114").Value = -65018
$ws.Range("H122").Value = 2939.8386
$ws.Range("I122").Value = 2251.963
$ws.Range("J122").Value = 7583
$ws.Range("K122").Value = 6755.889000000001
$ws.Range("L122").Value = 22749
$ws.Range("M122").Value = -4305.889000000001
$ws.Range("N122").Value = -27649

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6496612.5
$ws.Range("I99").Value = 2567.625
$ws.Range("J99").Value = 15155339
$ws.Range("K99").Value = 2567.625
$ws.Range("L99").Value = 15155339
$ws.Range("N99").Value = -15158335
$ws.Range("H105").Value = 3421.861
$ws.Range("I105").Value = 2617.261
$ws.Range("K105").Value = 2617.261
$ws.Range("M105").Value = -870.261
$ws.Range("H134").Value = 3971374.5
$ws.Range("I134").Value = 4809806
$ws.Range("K134").Value = 14429418
$ws.Range("M134").Value = -14426883

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4581.3716
$ws.Range("J132").Value = 5832.1113
$ws.Range("L132").Value = 17496.3339
$ws.Range("N132").Value = -22556.3339
$ws.Range("H134").Value = 3273.3333
$ws.Range("I134").Value = 1458.2963
$ws.Range("J134").Value = 5315.25
$ws.Range("K134").Value = 4374.8889
$ws.Range("L134").Value = 15945.75
$ws.Range("M134").Value = -1839.8889
$ws.Range("N134").Value = -21015.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 159.2
$ws.Range("I8").Value = 159.2
$ws.Range("K8").Value = 477.6
$ws.Range("M8").Value = -338.6
$ws.Range("H11").Value = 2491
$ws.Range("I11").Value = 477.5
$ws.Range("J11").Value = 3833.3333
$ws.Range("K11").Value = 1432.5
$ws.Range("L11").Value = 11499.9999
$ws.Range("M11").Value = -1292.5
$ws.Range("N11").Value = -11779.9999
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("H26").Value = 525.5
$ws.Range("J26").Value = 525.5
$ws.Range("L26").Value = 1576.5
$ws.Range("N26").Value = -2152.5
$ws.Range("H86").Value = 549
$ws.Range("I86").Value = 185.5
$ws.Range("J86").Value = 2003
$ws.Range("K86").Value = 556.5
$ws.Range("L86").Value = 6009
$ws.Range("M86").Value = 629.5
$ws.Range("N86").Value = -8381
$ws.Range("H89").Value = 549
$ws.Range("I89").Value = 185.5
$ws.Range("J89").Value = 2003
$ws.Range("K89").Value = 1669.5
$ws.Range("L89").Value = 18027
$ws.Range("M89").Value = 4258.5
$ws.Range("N89").Value = -29883
$ws.Range("H98").Value = 1318.25
$ws.Range("J98").Value = 1798.8
$ws.Range("L98").Value = 5396.4
$ws.Range("N98").Value = -8392.4
$ws.Range("H122").Value = 2358091.5
$ws.Range("I122").Value = 5657816.5
$ws.Range("K122").Value = 50920348.5
$ws.Range("M122").Value = -50917898.5
$ws.Range("H131").Value = 2125.5454
$ws.Range("I131").Value = 1194
$ws.Range("J131").Value = 3755.75
$ws.Range("K131").Value = 3582
$ws.Range("L131").Value = 11267.25
$ws.Range("M131").Value = 1458
$ws.Range("N131").Value = -21347.25
$ws.Range("H137").Value = 120598.06
$ws.Range("J137").Value = 128347.625
$ws.Range("L137").Value = 385042.875
$ws.Range("N137").Value = -395242.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 102267.5
$ws.Range("J80").Value = 202074.2
$ws.Range("L80").Value = 202074.2
$ws.Range("N80").Value = -204070.2
$ws.Range("H83").Value = 102267.5
$ws.Range("J83").Value = 202074.2
$ws.Range("L83").Value = 1010371
$ws.Range("N83").Value = -1020355
$ws.Range("H113").Value = 5589.64
$ws.Range("I113").Value = 3570.3225
$ws.Range("K113").Value = 3570.3225
$ws.Range("M113").Value = -1400.3225
$ws.Range("H122").Value = 9083071
$ws.Range("I122").Value = 14528512
$ws.Range("K122").Value = 43585536
$ws.Range("M122").Value = -43583086
$ws.Range("H126").Value = 4300.8125
$ws.Range("I126").Value = 1337.5
$ws.Range("K126").Value = 4012.5
$ws.Range("M126").Value = -1542.5
$ws.Range("H132").Value = 7249.5
$ws.Range("I132").Value = 2499.5
$ws.Range("J132").Value = 11999.5
$ws.Range("K132").Value = 7498.5
$ws.Range("L132").Value = 35998.5
$ws.Range("M132").Value = -4968.5
$ws.Range("N132").Value = -41058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2721.8572
$ws.Range("I22").Value = 1863.1818
$ws.Range("K22").Value = 1863.1818
$ws.Range("M22").Value = -1568.1818
$ws.Range("H27").Value = 2721.8572
$ws.Range("I27").Value = 1863.1818
$ws.Range("K27").Value = 1863.1818
$ws.Range("M27").Value = -1756.1818
$ws.Range("H136").Value = 8499.352999999999
$ws.Range("I136").Value = 2599.0334
$ws.Range("J136").Value = 16928.38
$ws.Range("K136").Value = 7797.100199999999
$ws.Range("L136").Value = 50785.14
$ws.Range("M136").Value = -5247.100199999999
$ws.Range("N136").Value = -55885.14

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 23350954
$ws.Range("I81").Value = 1429798.9
$ws.Range("K81").Value = 2859597.8
$ws.Range("M81").Value = -2858536.8
$ws.Range("H84").Value = 23350954
$ws.Range("I84").Value = 1429798.9
$ws.Range("K84").Value = 14297989
$ws.Range("M84").Value = -14292685
$ws.Range("H107").Value = 25642918
$ws.Range("I107").Value = 2283.3333
$ws.Range("K107").Value = 6849.999899999999
$ws.Range("M107").Value = -4929.999899999999
$ws.Range("H113").Value = 1132.8667
$ws.Range("I113").Value = 1020.34485
$ws.Range("J113").Value = 1336.8125
$ws.Range("K113").Value = 3061.03455
$ws.Range("L113").Value = 4010.4375
$ws.Range("M113").Value = -891.0345499999999
$ws.Range("N113").Value = -8350.4375
$ws.Range("H122").Value = 151588.89
$ws.Range("I122").Value = 192672.42
$ws.Range("K122").Value = 578017.26
$ws.Range("M122").Value = -575567.26
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0

# Special-case cell insertions/removals (cells added or removed entirely in the row)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M17").Value = -2832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M99").Value = -1069.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M20").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N126").ClearContents()
